$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 119 (existing rows 119.. shift down to 121..)
$ws.Rows.Item(119).Insert()
$ws.Rows.Item(119).Insert()

# Fill new row 119: Larry Ann / Primera
$ws.Range("A119").Value = 5
$ws.Range("B119").Value = "Macroferia Regional de Talca"
$ws.Range("C119").Value = "Maule"
$ws.Range("D119").Value = 44981
$ws.Range("E119").Value = 7
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100103
$ws.Range("H119").Value = "Frutos de hueso (carozo)"
$ws.Range("I119").Value = 100103002
$ws.Range("J119").Value = "Ciruela"
$ws.Range("K119").Value = "Larry Ann"
$ws.Range("L119").Value = "Primera"
$ws.Range("M119").Value = 210
$ws.Range("N119").Value = 10000
$ws.Range("O119").Value = 10000
$ws.Range("P119").Value = 10000
$ws.Range("Q119").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R119").Value = "Provincia de Curicó"
$ws.Range("S119").Value = 556
$ws.Range("T119").Value = 18

# Fill new row 120: Larry Ann / Segunda
$ws.Range("A120").Value = 5
$ws.Range("B120").Value = "Macroferia Regional de Talca"
$ws.Range("C120").Value = "Maule"
$ws.Range("D120").Value = 44981
$ws.Range("E120").Value = 7
$ws.Range("F120").Value = "Fruta"
$ws.Range("G120").Value = 100103
$ws.Range("H120").Value = "Frutos de hueso (carozo)"
$ws.Range("I120").Value = 100103002
$ws.Range("J120").Value = "Ciruela"
$ws.Range("K120").Value = "Larry Ann"
$ws.Range("L120").Value = "Segunda"
$ws.Range("M120").Value = 250
$ws.Range("N120").Value = 8000
$ws.Range("O120").Value = 8000
$ws.Range("P120").Value = 8000
$ws.Range("Q120").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R120").Value = "Provincia de Curicó"
$ws.Range("S120").Value = 444
$ws.Range("T120").Value = 18

Write-Host "Done"
